# "Generate Report for Handoff"
# Re-running the handoff report updates the "Latest Handoff Datetime" (column D)
# for every file that was just handed off, on both the zh-cn and de-de status
# sheets. Rows 7, 10, 11, 12, 13, 14, 15, 16 correspond to files handed off in
# this batch.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$handoffRows = @(7, 10, 11, 12, 13, 14, 15, 16)

foreach ($r in $handoffRows) {
    $zhcn.Range("D$r").Value = "2016-03-08 06:15:13"
}

foreach ($r in $handoffRows) {
    $dede.Range("D$r").Value = "2016-03-08 06:15:17"
}
